# Daily COVID-19 country data refresh ("Update countries & provincias Spain"):
#  - Spain (F5), Iran (row 10), Austria (row 17), Australia (row 23),
#    Romania (row 30) and Bosnia y Herzegovina (row 72) get new case counts.
#  - The country list is re-sorted by total cases (column B) descending, which
#    reshuffles a number of same-count rows.
#  - A new country, "Sudan del Sur", is appended as row 213.
#  - The "Datos actualizados" timestamp (A1) is refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Refresh the "datos actualizados" timestamp header (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 5 de Abril de 2020 a las 12:22"

# --- Full country table (rows 4-213): País, Casos totales, Nuevos casos,
#     Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes ---
$rows = @(
  @("Estados Unidos",311637,280,14828,288355,8206,2,8454),
  @("España",130759,4591,38080,80261,6861,471,12418),
  @("Italia",124632,0,20996,88274,3994,0,15362),
  @("Alemania",96108,16,26400,68264,3936,0,1444),
  @("Francia",89953,0,15438,66955,6838,0,7560),
  @("China",81669,30,76964,1376,295,3,3329),
  @("Iran",58226,2483,19736,34887,4103,151,3603),
  @("Reino Unido",41903,0,135,37455,1559,0,4313),
  @("Turquia",23934,0,786,22647,1311,0,501),
  @("Suiza",20510,5,6415,13429,391,0,666),
  @("Belgica",19691,1260,3751,14493,1261,164,1447),
  @("Paises Bajos",16627,0,250,14726,1360,0,1651),
  @("Canada",14018,106,2603,11182,426,2,233),
  @("Austria",11850,69,2998,8648,244,18,204),
  @("Portugal",10524,0,75,10183,251,0,266),
  @("Brasil",10360,0,127,9788,296,0,445),
  @("Corea del Sur",10237,81,6463,3591,55,6,183),
  @("Israel",8018,167,477,7495,127,2,46),
  @("Suecia",6443,0,205,5865,520,0,373),
  @("Australia",5687,137,2315,3338,85,4,34),
  @("Noruega",5645,95,32,5551,98,0,62),
  @("Rusia",4731,0,333,4355,8,0,43),
  @("Irlanda",4604,0,25,4442,148,0,137),
  @("Chequia",4475,3,78,4335,86,3,62),
  @("Chile",4161,0,528,3606,38,0,27),
  @("Dinamarca",4077,0,1283,2633,142,0,161),
  @("Rumania",3864,251,374,3342,141,2,148),
  @("Polonia",3834,207,134,3616,50,5,84),
  @("Malasia",3662,179,1005,2596,99,4,61),
  @("India",3588,0,229,3260,0,0,99),
  @("Ecuador",3465,0,100,3193,100,0,172),
  @("Filipinas",3246,152,64,3030,1,8,152),
  @("Japon",3139,0,514,2548,64,0,77),
  @("Pakistan",2880,62,170,2665,18,4,45),
  @("Luxemburgo",2729,0,500,2198,35,0,31),
  @("Arabia Saudita",2370,191,420,1921,41,0,29),
  @("Indonesia",2273,181,164,1911,0,7,198),
  @("Tailandia",2169,102,674,1472,23,3,23),
  @("Finlandia",1927,45,300,1602,73,0,25),
  @("Mexico",1890,202,633,1178,1,19,79),
  @("Panama",1801,0,13,1742,75,0,46),
  @("Peru",1746,0,914,759,88,0,73),
  @("Grecia",1673,0,78,1527,92,0,68),
  @("Serbia",1624,0,54,1526,89,0,44),
  @("Sudafrica",1585,0,95,1481,7,0,9),
  @("Republica Dominicana",1578,0,17,1484,147,0,77),
  @("Emiratos Arabes Unidos",1505,0,125,1370,2,0,10),
  @("Argentina",1451,0,279,1129,86,0,43),
  @("Islandia",1417,0,396,1017,12,0,4),
  @("Colombia",1406,0,85,1289,50,0,32),
  @("Catar",1325,0,109,1213,37,0,3),
  @("Ucrania",1251,26,25,1194,16,0,32),
  @("Argelia",1251,0,90,1031,0,0,130),
  @("Singapur",1189,0,297,886,24,0,6),
  @("Croacia",1126,0,119,995,39,0,12),
  @("Estonia",1097,58,62,1020,17,2,15),
  @("Egipto",1070,0,241,758,0,0,71),
  @("Nueva Zelanda",1039,89,156,882,3,0,1),
  @("Eslovenia",997,20,79,890,31,6,28),
  @("Marruecos",960,41,69,825,1,7,66),
  @("Hong Kong",890,28,206,680,8,0,4),
  @("Irak",878,0,259,563,0,0,56),
  @("Armenia",822,52,57,758,30,0,7),
  @("Lituania",811,40,7,793,11,0,11),
  @("Moldavia",752,0,30,708,65,2,14),
  @("Hungria",733,55,66,633,17,2,34),
  @("Crucero",712,0,619,82,10,0,11),
  @("Barein",688,0,423,261,3,0,4),
  @("Bosnia y Herzegovina",634,10,30,583,4,0,21),
  @("Kazajistan",569,38,36,527,6,1,6),
  @("Kuwait",556,77,99,456,17,0,1),
  @("Camerun",555,0,17,529,0,0,9),
  @("Tunez",553,0,5,530,26,0,18),
  @("Letonia",533,24,1,531,4,0,1),
  @("Libano",527,7,54,455,26,1,18),
  @("Bulgaria",522,19,37,467,26,1,18),
  @("Azerbaiyan",521,0,32,484,17,0,5),
  @("Eslovaquia",485,14,10,474,3,0,1),
  @("Republica de Macedonia",483,0,20,446,11,0,17),
  @("Principado de Andorra",466,0,21,428,12,0,17),
  @("Bielorrusia",440,0,53,382,11,0,5),
  @("Costa Rica",435,0,13,420,13,0,2),
  @("Republica de Chipre",426,0,33,384,11,0,9),
  @("Uruguay",400,0,93,302,13,0,5),
  @("Taiwan",363,8,54,304,0,0,5),
  @("Afganistan",337,28,15,315,0,0,7),
  @("Reunion",334,0,40,294,4,0,0),
  @("Albania",333,0,99,214,7,0,20),
  @("Jordania",323,0,74,244,5,0,5),
  @("Burkina Faso",318,0,66,236,0,0,16),
  @("Uzbekistan",298,32,25,271,8,0,2),
  @("Oman",298,21,61,235,3,0,2),
  @("Cuba",288,0,15,267,11,0,6),
  @("Honduras",268,4,6,240,10,7,22),
  @("San Marino",259,0,27,200,14,0,32),
  @("Costa de Marfil",245,0,25,219,0,0,1),
  @("Vietnam",240,0,90,150,3,0,0),
  @("Mauricio",227,31,7,213,1,0,7),
  @("Estado de Palestina",226,9,23,202,0,0,1),
  @("Senegal",219,0,72,145,1,0,2),
  @("Nigeria",214,0,25,185,2,0,4),
  @("Ghana",214,9,31,178,2,0,5),
  @("Malta",213,0,2,211,3,0,0),
  @("Montenegro",203,2,1,200,4,0,2),
  @("Islas Feroe",181,0,99,82,1,0,0),
  @("Georgia",170,8,36,132,6,1,2),
  @("Sri Lanka",167,1,29,133,5,0,5),
  @("Bolivia",157,18,2,145,3,0,10),
  @("Venezuela",155,0,52,96,6,0,7),
  @("Consejo Danes para los Refugiados",154,0,3,133,0,0,18),
  @("Kirguistan",147,3,9,137,5,0,1),
  @("Martinica",145,0,27,115,22,0,3),
  @("Niger",144,0,0,136,0,0,8),
  @("Brunei",135,0,73,61,3,0,1),
  @("Mayotte",134,0,14,118,3,0,2),
  @("Guadalupe",134,0,24,103,14,0,7),
  @("Isla de Man",127,1,0,126,0,0,1),
  @("Kenia",126,0,4,118,2,0,4),
  @("Camboya",114,0,50,64,1,0,0),
  @("Guinea",111,0,5,106,0,0,0),
  @("Paraguay",104,8,12,89,2,0,3),
  @("Trinidad yTobago",103,0,1,96,0,0,6),
  @("Ruanda",102,0,0,102,0,0,0),
  @("Gibraltar",98,0,52,46,0,0,0),
  @("Banglades",88,18,30,49,1,1,9),
  @("Liechtenstein",77,0,0,76,0,0,1),
  @("Madagascar",70,0,0,70,6,0,0),
  @("Monaco",66,0,3,62,2,0,1),
  @("Aruba",64,0,1,63,0,0,0),
  @("El Salvador",62,6,2,57,4,0,3),
  @("Guatemala",61,0,15,44,1,0,2),
  @("Guayana Francesa",61,0,22,39,1,0,0),
  @("Jamaica",55,2,7,45,0,0,3),
  @("Barbados",52,0,0,52,0,0,0),
  @("Republica de Yibuti",50,0,8,42,0,0,0),
  @("Uganda",48,0,0,48,0,0,0),
  @("Congo",45,23,2,38,0,3,5),
  @("Macao",44,0,10,34,0,0,0),
  @("Mali",41,0,1,37,0,0,3),
  @("Togo",41,0,17,21,0,0,3),
  @("Polinesia Francesa",40,0,0,40,1,0,0),
  @("Puerto Rico",39,0,1,36,0,0,2),
  @("Zambia",39,0,2,36,0,0,1),
  @("Etiopia",38,0,4,34,1,0,0),
  @("Bermudas",37,0,14,23,0,0,0),
  @("Islas Caimanes",35,0,1,33,0,0,1),
  @("Guam",32,0,0,31,0,0,1),
  @("Eritrea",29,0,0,29,0,0,0),
  @("San Martin (Parte Francesa)",29,5,7,20,6,0,2),
  @("Bahamas",28,0,0,24,1,0,4),
  @("San Martin (Parte Holandesa)",25,2,6,15,0,0,4),
  @("Guyana",24,1,0,20,0,0,4),
  @("Haiti",21,1,1,20,0,0,0),
  @("Birmania",21,0,0,20,0,0,1),
  @("Gabon",21,0,1,19,0,0,1),
  @("Tanzania",20,0,3,16,0,0,1),
  @("Siria",19,3,2,15,0,0,2),
  @("Maldivas",19,0,13,6,0,0,0),
  @("Guinea-Bisau",18,0,0,18,0,0,0),
  @("Libia",18,0,0,17,0,0,1),
  @("Islas Virgenes de los Estados Unidos",17,0,0,17,0,0,0),
  @("Nueva Caledonia",17,0,1,16,0,0,0),
  @("Guinea Ecuatorial",16,0,1,15,0,0,0),
  @("Benin",16,0,2,14,0,0,0),
  @("Namibia",16,2,3,13,0,0,0),
  @("Antigua y Barbuda",15,0,0,15,1,0,0),
  @("Dominica",14,0,0,14,0,0,0),
  @("Santa Lucia",14,0,1,13,0,0,0),
  @("Mongolia",14,0,2,12,0,0,0),
  @("Granada",12,0,0,12,2,0,0),
  @("Fiyi",12,0,0,12,0,0,0),
  @("Laos",11,1,0,11,0,0,0),
  @("Groenlandia",11,0,3,8,0,0,0),
  @("Curazao",11,0,5,5,0,0,1),
  @("Seychelles",10,0,0,10,0,0,0),
  @("Mozambique",10,0,1,9,0,0,0),
  @("Surinam",10,0,0,9,0,0,1),
  @("Angola",10,0,2,6,0,0,2),
  @("Sudan",10,0,2,6,0,0,2),
  @("Liberia",10,0,3,6,0,0,1),
  @("Republica del Chad",9,0,0,9,0,0,0),
  @("Suazilandia",9,0,0,9,0,0,0),
  @("San Cristobal y Nieves",9,0,0,9,0,0,0),
  @("Nepal",9,0,1,8,0,0,0),
  @("Zimbabue",9,0,0,8,0,0,1),
  @("Montserrat",9,0,0,7,0,0,2),
  @("Republica de Africa Central",8,0,0,8,0,0,0),
  @("Santa Sede",7,0,0,7,0,0,0),
  @("Cabo Verde",7,0,0,6,0,0,1),
  @("San Vicente y las Granadinas",7,0,1,6,0,0,0),
  @("Somalia",7,0,1,6,0,0,0),
  @("San Bartolome",6,0,1,5,0,0,0),
  @("Mauritania",6,0,2,3,0,0,1),
  @("Islas Turcas y Caicos",5,0,0,5,0,0,0),
  @("Nicaragua",5,0,0,4,0,0,1),
  @("Butan",5,0,2,3,0,0,0),
  @("Sahara Occidental",4,0,0,4,0,0,0),
  @("Sierra Leona",4,0,0,4,0,0,0),
  @("Malaui",4,0,0,4,0,0,0),
  @("Belice",4,0,0,4,0,0,0),
  @("Botsuana",4,0,0,3,0,0,1),
  @("Gambia",4,0,2,1,0,0,1),
  @("Islas Virgenes Britanicas",3,0,0,3,0,0,0),
  @("Anguila",3,0,0,3,0,0,0),
  @("Burundi",3,0,0,3,0,0,0),
  @("Bonaire, San Eustaquio y Saba",2,0,0,2,0,0,0),
  @("Islas Malvinas",1,0,0,1,0,0,0),
  @("Papua Nueva Guinea",1,0,0,1,0,0,0),
  @("Timor Oriental",1,0,0,1,0,0,0),
  @("Sudan del Sur",1,1,0,1,0,0,0)
)

$numRows = $rows.Count
$numCols = 8
$data = New-Object 'object[,]' $numRows,$numCols
for ($i = 0; $i -lt $numRows; $i++) {
  for ($j = 0; $j -lt $numCols; $j++) {
    $data[$i,$j] = $rows[$i][$j]
  }
}

$ws.Range("A4:H213").Value = $data
